$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need special handling so they
# remain stored as text (matching the original inline-string cell type) instead of
# being auto-converted into a numeric cell by Excel.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "39.649.75"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.185.59"
$ws.Range("E3").Value = "  -2.32%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "290.13"
$ws.Range("E5").Value = "  -1.13%  "
Set-TextValue $ws.Range("D6") "85.58"
$ws.Range("E6").Value = "  -1.98%  "
Set-TextValue $ws.Range("D7") "0.505"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("E8").Value = "  -0.09%  "
Set-TextValue $ws.Range("D9") "0.463"
$ws.Range("E9").Value = "  -2.71%  "
Set-TextValue $ws.Range("D10") "30.10"
$ws.Range("E10").Value = "  -3.93%  "
Set-TextValue $ws.Range("D11") "49.94"
$ws.Range("E11").Value = "  +6.26%  "
Set-TextValue $ws.Range("D12") "0.0775"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("E13").Value = "  +2.39%  "
Set-TextValue $ws.Range("D14") "6.39"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "2.518.38"
$ws.Range("E15").Value = "  -2.63%  "
Set-TextValue $ws.Range("D16") "13.62"
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").Value = "2.205.61"
$ws.Range("E17").Value = "  -1.52%  "
Set-TextValue $ws.Range("D18") "0.722"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "39.512.17"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "0.0₃0878"
$ws.Range("E20").Value = "  -1.48%  "
Set-TextValue $ws.Range("D21") "11.07"
$ws.Range("E21").Value = "  -1.83%  "
Set-TextValue $ws.Range("D22") "5.68"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("E23").Value = "  -1.36%  "
Set-TextValue $ws.Range("D24") "236.11"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +0.24%  "
Set-TextValue $ws.Range("D26") "2.43"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").Value = "  -3.36%  "
Set-TextValue $ws.Range("D28") "22.95"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  -3.73%  "
Set-TextValue $ws.Range("D30") "9.12"
$ws.Range("E30").Value = "  -2.44%  "
Set-TextValue $ws.Range("D31") "155.66"
$ws.Range("E31").Value = "  +2.67%  "
Set-TextValue $ws.Range("D32") "31.02"
$ws.Range("E32").Value = "  -7.30%  "
$ws.Range("E33").Value = "  +0.09%  "
Set-TextValue $ws.Range("D34") "4.89"
$ws.Range("E34").Value = "  -0.98%  "
Set-TextValue $ws.Range("D35") "0.0702"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -2.35%  "
Set-TextValue $ws.Range("D37") "2.82"
$ws.Range("E37").Value = "  -0.34%  "
Set-TextValue $ws.Range("D38") "0.112"
$ws.Range("E38").Value = "  -0.39%  "
Set-TextValue $ws.Range("D39") "0.0967"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D40") "15.02"
$ws.Range("E40").Value = "  -8.10%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D41") "1.66"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").Value = "2.109.07"
$ws.Range("E42").Value = "  +2.13%  "
Set-TextValue $ws.Range("D43") "3.70"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  -1.74%  "
Set-TextValue $ws.Range("D46") "9.66"
$ws.Range("E46").Value = "  -1.86%  "
Set-TextValue $ws.Range("D47") "17.07"
$ws.Range("E47").Value = "  -6.42%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "2.393.08"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -0.38%  "
Set-TextValue $ws.Range("D51") "87.53"
$ws.Range("E51").Value = "  -2.25%  "
